$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("K1").Value = "Parent Type"
$ws.Range("L1").Value = "Parent Id"
$ws.Range("K1:L1").Select()
